# TMC.xlsx edit: correct the "age" (J5) input and fix the K5 discount
# formula (it was referencing K2 instead of J5), then leave the cursor
# on J5 where the user was working.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet (and its formulas) is protected with the password "tmc"
# (see the sharedStrings hint: "mot de passe protection formules: tmc").
$ws.Unprotect("tmc")

# Correct the vehicle age used by the discount formula.
$ws.Range("J5").Value = 1

# Fix the bugged discount formula: it must depend on J5 (the age), not K2.
$ws.Range("K5").Formula = "=IF(J5<=5, 1-(0.1*J5), IF(J5<=15, 1-(0.1*5 + 0.05*(J5-5)), 1-(0.1*5 + 0.05*10)))"

# Re-protect the sheet with the same password.
$ws.Protect("tmc")

# Leave the active selection on J5, matching where the edit was made.
$null = $ws.Range("J5").Select()
